$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-6
# from 2023-09-05 (45174) to 2023-09-06 (45175)
$ws.Range("C2:C6").Value = 45175
